$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.189457535743713
$ws.Range("B1").Value = 1.761683821678162
$ws.Range("C1").Value = 6.555927753448486
$ws.Range("D1").Value = 2.285098075866699
$ws.Range("E1").Value = 1.192834496498108
